# Generate Report for handoff
# - the source markdown file was renamed 5ac20b70-...-md -> b4cd3aed-...-md
# - the xlf hash / filename changed, new handoff timestamps recorded
# - the previously-failed handoff row (088b3ce5-...-md / "Handoff transform failed")
#   is gone entirely (row deleted, not just cleared) on every sheet
$wb = $excel.ActiveWorkbook

$oldMd  = "5ac20b70-63ad-459b-9428-d49f65bd972c.md"
$newMd  = "b4cd3aed-69e7-4617-a156-447920c7b6c3.md"

$oldXlfZh = "5ac20b70-63ad-459b-9428-d49f65bd972c.cf4ccd8494ceb9bf3bb47a99591c2f5d31db45e6.zh-cn.xlf"
$newXlfZh = "b4cd3aed-69e7-4617-a156-447920c7b6c3.f0a5ec462cbd0913ce7fb53d9d42209889a21108.zh-cn.xlf"

$oldXlfDe = "5ac20b70-63ad-459b-9428-d49f65bd972c.cf4ccd8494ceb9bf3bb47a99591c2f5d31db45e6.de-de.xlf"
$newXlfDe = "b4cd3aed-69e7-4617-a156-447920c7b6c3.f0a5ec462cbd0913ce7fb53d9d42209889a21108.de-de.xlf"

$newZhTimestamp = "2016-02-16 15:20:44"
$newDeTimestamp = "2016-02-16 15:20:57"

$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/f699aeeff26bd23c1b0e558609177b5a09da0e1f/e2e/$newMd"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/f699aeeff26bd23c1b0e558609177b5a09da0e1f/.localization-config"
$xlfZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e00a5033af0b6bf6c79906785810330523e537f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newXlfZh"
$xlfDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4589c352b4cca285f246e01d4fd3bd0cbd260dbb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newXlfDe"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value2 = $newMd

# Drop the "Handoff transform failed" row (old row 3); row 4 shifts up to row 3.
$ws.Rows(3).Delete()

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value2 = $newMd
$ws.Range("C2").Value2 = $newXlfZh
$ws.Range("D2").Value2 = $newZhTimestamp

$ws.Rows(3).Delete()

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), $xlfZhUrl, "", "", $newXlfZh) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value2 = $newMd
$ws.Range("C2").Value2 = $newXlfDe
$ws.Range("D2").Value2 = $newDeTimestamp

$ws.Rows(3).Delete()

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), $xlfDeUrl, "", "", $newXlfDe) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $cfgUrl, "", "", ".localization-config") | Out-Null
